$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert a new column before column C ("default" values column), shifting
# the existing Min value / Max value / etc. columns one place to the right.
$ws.Columns("C:C").Insert()

# Try to align the new column's width/style with column B (best effort;
# the simulated ColumnWidth property only keeps 2 decimal digits of
# precision so this may not match the original raw width exactly).
$ws.Columns("C:C").ColumnWidth = $ws.Columns("B:B").ColumnWidth

# Header for the new column
$ws.Range("C1").Value = "default"

# Default values for each primitive type (new column C)
$ws.Range("C3").Value = 0      # byte
$ws.Range("C4").Value = 0      # short
$ws.Range("C5").Value = 0      # int
$ws.Range("C6").Value = 0      # long
$ws.Range("C7").Value = 0      # float
$ws.Range("C8").Value = 0      # double
$ws.Range("C9").Value = $false # boolean (stored as boolean type)

# Row 7 (float) previously had no styled cells past column A, so the new
# column C cell there doesn't inherit the shared centered style used by
# the other rows; apply it explicitly (xlCenter = -4108) to match.
$ws.Range("C7").HorizontalAlignment = -4108

# boolean row: replace "true/false" with "depends on JVM"
$ws.Range("B9").Value = "depends on JVM"

# New "Size" value for float row (row 7), which previously had no size.
$ws.Range("B7").Value = "32bit"

# The char row (10) has no "default" column value; the Insert() operation
# leaves a blank formatted placeholder behind in the new column, so remove
# it completely (no cell at all in that position).
$ws.Range("C10").Clear()

# Update selection to match the authored workbook
$ws.Range("D8").Select()
